$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "ML model retrained with all data": the retrained model now weights the
# full dataset (J) at 1 and the held-out slice (K) at 0.3 for every row,
# instead of only the first row being flagged ("r"/"s") while the rest of
# the rows used a 0.3/1 split. Swap/normalize columns J (10) and K (11)
# across every data row so J=1 and K=0.3 throughout.
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 10).Value = 1
    $ws.Cells.Item($r, 11).Value = 0.3
}

# Reflect the reviewer's scrolled/selected view: column K selected top to
# bottom, with the window scrolled so row 26 is at the top of the pane.
$ws.Range("K1:K$lastRow").Select()
$excel.ActiveWindow.ScrollRow = 26
$excel.ActiveWindow.ScrollColumn = 1
